$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells: force text interpretation so numeric-looking
# strings (e.g. "1.00", "525.43") are stored as text, not auto-converted
# to numbers, matching the original inlineStr cell type. The NumberFormat
# "@" forces text entry; resetting Style to "Normal" afterwards drops the
# temporary number-format style so no stray style id is left on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '64.141.90'
$ws.Range('E2').Value = '  +1.89%  '
Set-TextValue $ws.Range('D3') '3.385.69'
$ws.Range('E3').Value = '  +4.14%  '
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.05%  '
Set-TextValue $ws.Range('D5') '525.43'
$ws.Range('E5').Value = '  +2.21%  '
Set-TextValue $ws.Range('D6') '175.29'
$ws.Range('E6').Value = '  -1.15%  '
Set-TextValue $ws.Range('D7') '0.596'
$ws.Range('E7').Value = '  +0.11%  '
Set-TextValue $ws.Range('D8') '3.384.57'
$ws.Range('E8').Value = '  +4.18%  '
$ws.Range('E9').Value = '  -0.05%  '
Set-TextValue $ws.Range('D10') '0.611'
$ws.Range('E10').Value = '  -0.79%  '
Set-TextValue $ws.Range('D11') '53.69'
$ws.Range('E11').Value = '  -6.68%  '
Set-TextValue $ws.Range('D12') '0.135'
$ws.Range('E12').Value = '  +3.17%  '
$ws.Range('E13').Value = '  +1.61%  '
Set-TextValue $ws.Range('D14') '9.12'
$ws.Range('E14').Value = '  +0.55%  '
Set-TextValue $ws.Range('D15') '3.921.38'
$ws.Range('E15').Value = '  +3.83%  '
Set-TextValue $ws.Range('D16') '3.381.41'
$ws.Range('E16').Value = '  +3.89%  '
$ws.Range('E17').Value = '  +1.38%  '
Set-TextValue $ws.Range('D18') '17.64'
$ws.Range('E18').Value = '  +0.96%  '
Set-TextValue $ws.Range('D19') '64.101.88'
$ws.Range('E19').Value = '  +2.06%  '
Set-TextValue $ws.Range('D20') '11.29'
$ws.Range('E20').Value = '  +3.55%  '
$ws.Range('E21').Value = '  +2.28%  '
Set-TextValue $ws.Range('D22') '375.84'
$ws.Range('E22').Value = '  +1.48%  '
Set-TextValue $ws.Range('D23') '11.60'
$ws.Range('E23').Value = '  +3.47%  '
Set-TextValue $ws.Range('D24') '4.19'
$ws.Range('E24').Value = '  +9.53%  '
Set-TextValue $ws.Range('D25') '81.58'
$ws.Range('E25').Value = '  +2.26%  '
Set-TextValue $ws.Range('D26') '3.70'
$ws.Range('E26').Value = '  +0.29%  '
Set-TextValue $ws.Range('D27') '6.16'
$ws.Range('E27').Value = '  +1.58%  '
Set-TextValue $ws.Range('D28') '2.72'
$ws.Range('E28').Value = '  +3.76%  '
Set-TextValue $ws.Range('D29') '11.36'
$ws.Range('E29').Value = '  -0.09%  '
Set-TextValue $ws.Range('D30') '8.26'
$ws.Range('E30').Value = '  -0.33%  '
Set-TextValue $ws.Range('D31') '29.09'
$ws.Range('E31').Value = '  +2.57%  '
Set-TextValue $ws.Range('D32') '631.96'
$ws.Range('E32').Value = '  -0.07%  '
Set-TextValue $ws.Range('D33') '6.48'
$ws.Range('E33').Value = '  -4.77%  '
Set-TextValue $ws.Range('D34') '11.26'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  +0.34%  '
Set-TextValue $ws.Range('D36') '57.99'
$ws.Range('E36').Value = '  -1.19%  '
$ws.Range('E37').Value = '  +0.03%  '
Set-TextValue $ws.Range('D38') '36.43'
$ws.Range('E38').Value = '  +0.46%  '
Set-TextValue $ws.Range('D39') '0.382'
$ws.Range('E39').Value = '  -4.70%  '
Set-TextValue $ws.Range('D40') '0.0₃0743'
$ws.Range('E40').Value = '  +13.29%  '
Set-TextValue $ws.Range('D41') '1.00'
$ws.Range('E41').Value = '  +0.25%  '
Set-TextValue $ws.Range('D42') '2.71'
$ws.Range('E42').Value = '  +10.47%  '
Set-TextValue $ws.Range('D43') '2.972.85'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('E44').Value = '  +1.27%  '
Set-TextValue $ws.Range('D45') '3.02'
$ws.Range('E45').Value = '  +7.23%  '
Set-TextValue $ws.Range('D46') '2.67'
$ws.Range('E46').Value = '  +3.68%  '
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('E48').Value = '  -2.14%  '
Set-TextValue $ws.Range('D49') '3.01'
$ws.Range('E49').Value = '  +1.40%  '
Set-TextValue $ws.Range('D50') '0.125'
$ws.Range('E50').Value = '  +0.25%  '
Set-TextValue $ws.Range('D51') '136.96'
$ws.Range('E51').Value = '  +4.70%  '
